$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 161; this shifts the existing rows 161-184 down to 162-185
$ws.Rows.Item(161).Insert()

# Populate the new row 161 with the new weekly record
$ws.Range("A161").Value = 8
$ws.Range("B161").Value = "Terminal La Palmera de La Serena"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 45180
$ws.Range("D161").NumberFormat = $ws.Range("D162").NumberFormat
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 100114007
$ws.Range("G161").Value = "Jengibre"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 400
$ws.Range("K161").Value = 16500
$ws.Range("L161").Value = 17000
$ws.Range("M161").Value = 16750
$ws.Range("N161").Value = "$/caja 13 kilos"
$ws.Range("O161").Value = "Perú"
$ws.Range("P161").Value = 1288
$ws.Range("Q161").Value = 13
$ws.Range("R161").Value = "Hortaliza"
